$wb = $excel.ActiveWorkbook

# Rename the "General" sheet to "Table"
$wsTable = $wb.Worksheets.Item("General")
$wsTable.Name = "Table"

$wsCodelists = $wb.Worksheets.Item("Codelists")

# Move the selection/active cell on the Codelists sheet, and make the
# renamed Table sheet the active (selected) tab when the file is saved.
$wsCodelists.Activate()
$wsCodelists.Range("O68").Select()

$wsTable.Activate()
$wsTable.Range("B52").Select()
